# Remove the "syntok" dependency row from the Acknowledgments sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column A holds the dependency name; find the row whose name is "syntok".
$found = $ws.Columns.Item(1).Find("syntok")

if ($found -ne $null) {
    $row = $found.Row
    $ws.Rows.Item($row).Delete()
}
